# Swap the Code/Sale Rate/Qty/Value figures (columns B, E, F, G) between
# the two rows of each duplicate-item pair. The item description (C) and
# cost rate (D) remain the same for both rows in a pair; only B, E, F, G
# need to trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowValues($row1, $row2) {
    foreach ($col in @("B", "E", "F", "G")) {
        $cell1 = $ws.Range("$col$row1")
        $cell2 = $ws.Range("$col$row2")
        $val1 = $cell1.Value2
        $val2 = $cell2.Value2
        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}

$rowPairs = @(
    @(192, 193),
    @(219, 220),
    @(227, 228),
    @(232, 233),
    @(243, 244),
    @(364, 365),
    @(366, 367),
    @(372, 373),
    @(375, 376),
    @(380, 381),
    @(442, 443),
    @(572, 573)
)

foreach ($pair in $rowPairs) {
    Swap-RowValues $pair[0] $pair[1]
}

"done"
